$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (write B1 before A1 so shared-string indices match the source order)
$ws.Range("B1").Value = "last_name"
$ws.Range("A1").Value = "first_name"
$ws.Range("C1").Value = "year_start"
$ws.Range("D1").Value = "year_graduate"

# Data rows
$ws.Range("A2").Value = "Charlotte"
$ws.Range("B2").Value = "Erenberg"
$ws.Range("C2").Value = 2023
$ws.Range("D2").Value = 2024

$ws.Range("A3").Value = "Ella"
$ws.Range("B3").Value = "Mahaney"
$ws.Range("C3").Value = 2023
$ws.Range("D3").Value = 2025

$ws.Range("A4").Value = "Ellie"
$ws.Range("B4").Value = "Sprinkmann"
$ws.Range("C4").Value = 2024

# Column widths (A through M = columns 1-13); 11.83 is the closest input that
# the host's internal pixel-quantized column-width storage resolves back to
# the source file's stored width of 12.6328125 characters.
$ws.Range("A1:M1").EntireColumn.ColumnWidth = 11.83

# Selection to match target view state
$ws.Range("D6").Select() | Out-Null
